# Insert two new paragraphs at the very top of the document:
#   1) "CASOS DE USO ___...___" as an underlined heading line.
#   2) An empty paragraph (spacer) below it.
# Both precede whatever used to be the first paragraph in the body.

$d = $word.ActiveDocument

# Collapsed range at the very start of the document (before the existing
# first paragraph). Inserting two paragraph breaks here pushes the original
# content down and leaves two brand-new empty paragraphs in front of it.
$r = $d.Range(0, 0)
$r.InsertParagraphBefore()
$r.InsertParagraphBefore()

# The first paragraph in the document is now the first of the two new ones;
# give it the heading text and underline formatting.
$heading = $d.Paragraphs(1).Range
$heading.Text = "CASOS DE USO ___________________________________________________________"
$heading.Font.Underline = 1

# The second paragraph stays empty, acting as a spacer before the rest of
# the original document content.
